# Generate Report for Archive
#
# The localization-status report is regenerated and the row for
# "407bfbc6-984d-45e5-9b9d-7de10de0e775" moves from its old position
# (6th data row) up to become the 2nd data row, right after
# "b7fcfb28-3ab7-434f-9dfc-8ebe26d9fb78" on every sheet (Overview,
# zh-cn, de-de). All the other file rows keep their own data but shift
# down by one row to make room. This script captures each affected
# row's cell values + hyperlinks, then rewrites the block in the new
# row order so every sheet ends up matching the regenerated report.

$wb = $excel.ActiveWorkbook

# Desired new ordering of the file-id keys for the reordered block.
# (Row 1 is the header, the first data row "b7fcfb28..." and the last
# two data rows "95be4f59...","d99caf68..." plus the trailing
# ".localization-config" row do not move.)
$newKeyOrder = @(
    "407bfbc6-984d-45e5-9b9d-7de10de0e775",
    "576282a1-063c-44cd-bd5b-227707c7a3c9",
    "7d240ce6-9811-446e-ba47-5902a918d37c",
    "b96d7e22-f461-48af-95b0-d05193254f94",
    "e834abc6-3bbc-4ca8-b6b5-92e40f9ed347"
)

$firstRow = 3   # first row of the block being reordered
$lastRow  = 7   # last row of the block being reordered
$lastCol  = 9   # column I is the widest column in use (zh-cn / de-de sheets)

function Get-RowKey($ws, $row) {
    $v = $ws.Cells.Item($row, 1).Value()
    if ($v -eq $null) { return "" }
    $parts = $v.ToString().Split(".")
    return $parts[0]
}

foreach ($ws in $wb.Worksheets) {

    # ---- 1. snapshot existing hyperlinks on this sheet, keyed by "row,col"
    $hlMap = @{}
    foreach ($hl in $ws.Hyperlinks) {
        $r = $hl.Range.Row
        $c = $hl.Range.Column
        $key = "$r,$c"
        $hlMap[$key] = @{ addr = $hl.Address; disp = $hl.TextToDisplay }
    }

    # ---- 2. snapshot the rows that are about to be reordered
    $bundles = @{}
    for ($row = $firstRow; $row -le $lastRow; $row++) {
        $key = Get-RowKey $ws $row
        if ($key -eq "") { continue }

        $cells = @{}
        for ($col = 1; $col -le $lastCol; $col++) {
            $val = $ws.Cells.Item($row, $col).Value()
            $hlKey = "$row,$col"
            $hl = $null
            if ($hlMap.ContainsKey($hlKey)) { $hl = $hlMap[$hlKey] }
            $cells[$col] = @{ value = $val; hyperlink = $hl }
        }
        $bundles[$key] = $cells
    }

    if ($bundles.Count -eq 0) { continue }

    # ---- 3. clear the block (values + hyperlinks) before rewriting it
    for ($row = $firstRow; $row -le $lastRow; $row++) {
        for ($col = 1; $col -le $lastCol; $col++) {
            $cell = $ws.Cells.Item($row, $col)
            $cell.ClearContents()
        }
    }
    foreach ($hl in @($ws.Hyperlinks)) {
        $r = $hl.Range.Row
        if ($r -ge $firstRow -and $r -le $lastRow) {
            $hl.Delete()
        }
    }

    # ---- 4. write the bundles back out in the new order
    $destRow = $firstRow
    foreach ($key in $newKeyOrder) {
        if (-not $bundles.ContainsKey($key)) { $destRow++; continue }
        $cells = $bundles[$key]
        foreach ($col in $cells.Keys) {
            $data = $cells[$col]
            $cell = $ws.Cells.Item($destRow, $col)
            if ($data.value -ne $null) {
                $cell.Value = $data.value
            }
            if ($data.hyperlink -ne $null) {
                $ws.Hyperlinks.Add($cell, $data.hyperlink.addr, "", "", $data.hyperlink.disp) | Out-Null
            }
        }
        $destRow++
    }
}
